$wb = $excel.ActiveWorkbook

# ---- Sheet: 保險 (insurance) ----
$ws7 = $wb.Worksheets.Item("保險")
$ws7.Range("B1").Value = "company"
$ws7.Range("C1").Value = "name"
$ws7.Range("D1").Value = "owner"
$ws7.Range("E1").Value = "property_category"
$ws7.Range("F1").Value = "category"
$ws7.Range("G1").Value = "date"
$ws7.Range("H1").Value = "legislator_name"
$ws7.Range("I1").Value = "legislator_id"
$ws7.Range("J1").Value = "source_file"
$ws7.Range("K1").Value = "index"

$s7data = @(
    @(182,"南山人壽保險股份有限公司","南山康寧終身壽險","尤美女","insurance","normal","2013-12-20","尤美女",1730,"tmp84bd1",182),
    @(183,"南山人壽保險股份有限公司","南山人壽伴我一生躉繳變額壽險","尤美女","insurance","normal","2013-12-20","尤美女",1730,"tmp84bd1",183),
    @(184,"南山人壽保險股份有限公司","南山全新增額養老壽險","尤美女","insurance","normal","2013-12-20","尤美女",1730,"tmp84bd1",184),
    @(185,"南山人壽保險股份有限公司","南山金美滿還本終身保險","尤美女","insurance","normal","2013-12-20","尤美女",1730,"tmp84bd1",185),
    @(186,"南山人壽保險股份有限公司","南山美滿還本終身保險","尤美女","insurance","normal","2013-12-20","尤美女",1730,"tmp84bd1",186),
    @(187,"南山人壽保險股份有限公司","南山人壽新康祥終身保險一B型","尤美女","insurance","normal","2013-12-20","尤美女",1730,"tmp84bd1",187),
    @(188,"南山人壽保險股份有限公司","南山全新增額養老險壽","尤美女","insurance","normal","2013-12-20","尤美女",1730,"tmp84bd1",188),
    @(189,"南山人壽保險股份有限公司","南山金美滿還本終身保險","尤美女","insurance","normal","2013-12-20","尤美女",1730,"tmp84bd1",189),
    @(190,"南山人壽保險股份有限公司","南山美滿還本終身保險","尤美女","insurance","normal","2013-12-20","尤美女",1730,"tmp84bd1",190),
    @(191,"南山人壽保險股份有限公司","南山人壽新康祥終身保險一B型","尤美女","insurance","normal","2013-12-20","尤美女",1730,"tmp84bd1",191),
    @(192,"南山人壽保險股份有限公司","南山全新增額養老壽險","尤美女","insurance","normal","2013-12-20","尤美女",1730,"tmp84bd1",192),
    @(193,"全球人壽","全球人壽樂活外幣變額年金保險","黃瑞明","insurance","normal","2013-12-20","尤美女",1730,"tmp84bd1",193),
    @(194,"安聯人壽","超優勢變額萬能壽險（丙型）","黃瑞明","insurance","normal","2013-12-20","尤美女",1730,"tmp84bd1",194),
    @(196,"安聯人壽","超優勢變額年金保險","黃瑞明","insurance","normal","2013-12-20","尤美女",1730,"tmp84bd1",196),
    @(197,"法國巴黎人壽","法國巴黎人壽富足人生變額年金保險","黃瑞明","insurance","normal","2013-12-20","尤美女",1730,"tmp84bd1",197),
    @(198,"富邦人壽","年金保險:心得意利變型年金","黃瑞明","insurance","normal","2013-12-20","尤美女",1730,"tmp84bd1",198),
    @(199,"第一金人壽","外幣保險:錢進富林","黃瑞明","insurance","normal","2013-12-20","尤美女",1730,"tmp84bd1",199),
    @(200,"第一金人壽","白金人壽保險:再接債利","黃瑞明","insurance","normal","2013-12-20","尤美女",1730,"tmp84bd1",200),
    @(201,"中華郵政","郵政簡易人壽六年期吉利保險","黃瑞明","insurance","normal","2013-12-20","尤美女",1730,"tmp84bd1",201),
    @(202,"法國巴黎人壽","致勝100(年金）","黃瑞明","insurance","normal","2013-12-20","尤美女",1730,"tmp84bd1",202),
    @(203,"全球人壽","好康再現專案","黃瑞明","insurance","normal","2013-12-20","尤美女",1730,"tmp84bd1",203),
    @(204,"安聯人壽","匯豐人壽新好醫日額醫療養老保險","尤美女","insurance","normal","2013-12-20","尤美女",1730,"tmp84bd1",204),
    @(205,"安聯人壽","匯豐人壽新大安心傷害保本保險","尤美女","insurance","normal","2013-12-20","尤美女",1730,"tmp84bd1",205),
    @(206,"南山人壽","南山金美滿還本終身保險","尤美女","insurance","normal","2013-12-20","尤美女",1730,"tmp84bd1",206),
    @(207,"南山人壽","南山美滿還本終身保險","尤美女","insurance","normal","2013-12-20","尤美女",1730,"tmp84bd1",207),
    @(208,"南山人壽","南山人壽新康祥終身壽險B型","尤美女","insurance","normal","2013-12-20","尤美女",1730,"tmp84bd1",208),
    @(209,"富邦人壽","美利成增","黃瑞明","insurance","normal","2013-12-20","尤美女",1730,"tmp84bd1",209),
    @(210,"第一金人壽","富貴臨門","黃瑞明","insurance","normal","2013-12-20","尤美女",1730,"tmp84bd1",210),
    @(211,"中國人壽","富美滿利變型年金","黃瑞明","insurance","normal","2013-12-20","尤美女",1730,"tmp84bd1",211),
    @(212,"富邦人壽","鑽美利增額壽險","黃瑞明","insurance","normal","2013-12-20","尤美女",1730,"tmp84bd1",212),
    @(213,"中泰人壽","新投資贏家（年金）","黃瑞明","insurance","normal","2013-12-20","尤美女",1730,"tmp84bd1",213),
    @(214,"南山人壽保險股份有限公司","南山康樂限期繳費終身壽險","黃瑞明","insurance","normal","2013-12-20","尤美女",1730,"tmp84bd1",214),
    @(215,"南山人壽保險股份有限公司","南山康福二十年期繳費終身壽險","黃瑞明","insurance","normal","2013-12-20","尤美女",1730,"tmp84bd1",215),
    @(216,"南山人壽保險股份有限公司","南山康寧終身壽險","黃瑞明","insurance","normal","2013-12-20","尤美女",1730,"tmp84bd1",216)
)

for ($i = 0; $i -lt $s7data.Count; $i++) {
    $r = $i + 2
    $row = $s7data[$i]
    $ws7.Cells.Item($r, 1).Value = $row[0]
    $ws7.Cells.Item($r, 2).Value = $row[1]
    $ws7.Cells.Item($r, 3).Value = $row[2]
    $ws7.Cells.Item($r, 4).Value = $row[3]
    $ws7.Cells.Item($r, 5).Value = $row[4]
    $ws7.Cells.Item($r, 6).Value = $row[5]
    $ws7.Cells.Item($r, 7).Value = $row[6]
    $ws7.Cells.Item($r, 8).Value = $row[7]
    $ws7.Cells.Item($r, 9).Value = $row[8]
    $ws7.Cells.Item($r, 10).Value = $row[9]
    $ws7.Cells.Item($r, 11).Value = $row[10]
}

# ---- Sheet: 事業投資 (investment) ----
$ws8 = $wb.Worksheets.Item("事業投資")
$ws8.Range("B1").Value = "owner"
$ws8.Range("C1").Value = "company"
$ws8.Range("D1").Value = "address"
$ws8.Range("E1").Value = "total"
$ws8.Range("F1").Value = "register_date"
$ws8.Range("G1").Value = "register_reason"
$ws8.Range("H1").Value = "property_category"
$ws8.Range("I1").Value = "category"
$ws8.Range("J1").Value = "date"
$ws8.Range("K1").Value = "legislator_name"
$ws8.Range("L1").Value = "legislator_id"
$ws8.Range("M1").Value = "source_file"
$ws8.Range("N1").Value = "index"

$s8data = @(
    @(229,"黃瑞明","先驅媒體社會企業股份有限公司","臺北市仁愛路2段98號7樓",500000,"98年07月02日","投資","investment","normal","2013-12-20","尤美女",1730,"tmp84bd1",229),
    @(230,"黃瑞明","能得科技股份有限公司","臺北市敦化南路1段132號5樓",100000,"95年01月01日","投資","investment","normal","2013-12-20","尤美女",1730,"tmp84bd1",230)
)

for ($i = 0; $i -lt $s8data.Count; $i++) {
    $r = $i + 2
    $row = $s8data[$i]
    $ws8.Cells.Item($r, 1).Value = $row[0]
    $ws8.Cells.Item($r, 2).Value = $row[1]
    $ws8.Cells.Item($r, 3).Value = $row[2]
    $ws8.Cells.Item($r, 4).Value = $row[3]
    $ws8.Cells.Item($r, 5).Value = $row[4]
    $ws8.Cells.Item($r, 6).Value = $row[5]
    $ws8.Cells.Item($r, 7).Value = $row[6]
    $ws8.Cells.Item($r, 8).Value = $row[7]
    $ws8.Cells.Item($r, 9).Value = $row[8]
    $ws8.Cells.Item($r, 10).Value = $row[9]
    $ws8.Cells.Item($r, 11).Value = $row[10]
    $ws8.Cells.Item($r, 12).Value = $row[11]
    $ws8.Cells.Item($r, 13).Value = $row[12]
    $ws8.Cells.Item($r, 14).Value = $row[13]
}
